$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 165 ----
$ws.Range("A157").Copy()
$ws.Range("A165").PasteSpecial(-4122)
$ws.Range("A165").Value = "Model 2.4"
$ws.Range("B157").Copy()
$ws.Range("B165").PasteSpecial(-4122)
$ws.Range("B165").Value = "(64, 64, 1)"
$ws.Range("C165").Value = 32
$ws.Range("I165").Value = 0.0005
$ws.Range("T165").Formula = "=L165-J165"

# ---- Row 166 ----
$ws.Range("A158").Copy()
$ws.Range("A166").PasteSpecial(-4122)
$ws.Range("A166").Value = "Augmentation"
$ws.Range("O166").Value = 0.2
$ws.Range("P166").Value = 0.2
$ws.Range("Q166").Value = "(0,2, 0,2)"
$ws.Range("R166").Value = "hor"

# ---- Row 167 ----
$ws.Range("A167").Value = "Conv 1"
$ws.Range("D167").Value = 16
$ws.Range("E159").Copy()
$ws.Range("E167").PasteSpecial(-4122)
$ws.Range("E167").Value = "(3,3)"
$ws.Range("F159").Copy()
$ws.Range("F167").PasteSpecial(-4122)
$ws.Range("F167").Value = "(2,2)"
$ws.Range("G159").Copy()
$ws.Range("G167").PasteSpecial(-4122)
$ws.Range("G167").Value = "no"
$ws.Range("H167").Value = 0.2
$ws.Range("T167").Formula = "=L167-J167"

# ---- Row 168 ----
$ws.Range("A168").Value = "Conv 2"
$ws.Range("D168").Value = 64
$ws.Range("E160").Copy()
$ws.Range("E168").PasteSpecial(-4122)
$ws.Range("E168").Value = "(3,3)"
$ws.Range("F160").Copy()
$ws.Range("F168").PasteSpecial(-4122)
$ws.Range("F168").Value = "(2,2)"
$ws.Range("G160").Copy()
$ws.Range("G168").PasteSpecial(-4122)
$ws.Range("G168").Value = "no"
$ws.Range("H168").Value = 0.3
$ws.Range("T168").Formula = "=L168-J168"

# ---- Row 169 ----
$ws.Range("A169").Value = "Conv 3"
$ws.Range("D169").Value = 64
$ws.Range("E161").Copy()
$ws.Range("E169").PasteSpecial(-4122)
$ws.Range("E169").Value = "(5,5)"
$ws.Range("F161").Copy()
$ws.Range("F169").PasteSpecial(-4122)
$ws.Range("F169").Value = "(2,2)"
$ws.Range("G161").Copy()
$ws.Range("G169").PasteSpecial(-4122)
$ws.Range("G169").Value = "no"
$ws.Range("H169").Value = 0.5
$ws.Range("T169").Formula = "=L169-J169"

# ---- Row 170 ----
$ws.Range("A170").Value = "Dense"
$ws.Range("D170").Value = 128
$ws.Range("H170").Value = 0.5
$ws.Range("T170").Formula = "=L170-J170"

# ---- Row 171 ----
$ws.Range("A163").Copy()
$ws.Range("A171").PasteSpecial(-4122)
$ws.Range("A171").Value = "Output"
$ws.Range("B163").Copy()
$ws.Range("B171").PasteSpecial(-4122)
$ws.Range("C163").Copy()
$ws.Range("C171").PasteSpecial(-4122)
$ws.Range("D163").Copy()
$ws.Range("D171").PasteSpecial(-4122)
$ws.Range("E163").Copy()
$ws.Range("E171").PasteSpecial(-4122)
$ws.Range("F163").Copy()
$ws.Range("F171").PasteSpecial(-4122)
$ws.Range("G163").Copy()
$ws.Range("G171").PasteSpecial(-4122)
$ws.Range("H163").Copy()
$ws.Range("H171").PasteSpecial(-4122)
$ws.Range("I163").Copy()
$ws.Range("I171").PasteSpecial(-4122)
$ws.Range("J163").Copy()
$ws.Range("J171").PasteSpecial(-4122)
$ws.Range("J171").Value = 0.45490000000000003
$ws.Range("K163").Copy()
$ws.Range("K171").PasteSpecial(-4122)
$ws.Range("K171").Value = 0.79630000000000001
$ws.Range("L163").Copy()
$ws.Range("L171").PasteSpecial(-4122)
$ws.Range("L171").Value = 0.7742
$ws.Range("M163").Copy()
$ws.Range("M171").PasteSpecial(-4122)
$ws.Range("M171").Value = 0.64100000000000001
$ws.Range("N163").Copy()
$ws.Range("N171").PasteSpecial(-4122)
$ws.Range("N171").Value = 33
$ws.Range("T171").Formula = "=L171-J171"

# ---- Row height for the new thick-bottom summary row (row 171) ----
$ws.Rows.Item(171).RowHeight = 15

# ---- View state: freeze panes still at row 1, but scrolled so the new
#      block is visible, and the active selection matches the authored file ----
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("V171").Select()
